$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# Simplify the runs around "ou seja, cada imagem possui ..." - this text
# was originally split across several runs separated by grammar proofErr
# marks (", com 1,48MB... ou seja, " / "cada" / " " / "imagem" / " possui
# ... obtida"). The edit merges them and drops the proofErr marks while
# keeping the surrounding text untouched.
$r1 = $d.Content
$found1 = $r1.Find.Execute(", com 1,48MB (1.556.480 bytes em disco), ou seja, cada imagem possui exatamente metade do tamanho da original, nenhuma compressão obtida", $true, $false, $false, $false, $false, $true, 1, $false, ", com 1,48MB (1.556.480 bytes em disco), ou seja, cada imagem possui exatamente metade do tamanho da original, nenhuma compressão obtida", 2)

# --- Change 2 -----------------------------------------------------------
# Append a new sentence about the above/below strategy to the paragraph
# that explains the side-by-side strategy, and relocate the "_GoBack"
# bookmark so that it sits right after the newly appended sentence
# (instead of at the very end of the document).
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$r2 = $d.Content
$found2 = $r2.Find.Execute("o canal vermelho da imagem anaglífica.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Collapse(0)
$r2.InsertAfter(" A estratégia é similar para imagens do tipo above/below, sendo tratado a metade de cima da altura do container criado e depois a metade de baixo.")
$r2.Collapse(0)

$d.Bookmarks.Add("_GoBack", $r2)
